$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8927
$ws.Range("I116").Value = 21561
$ws.Range("J116").Value = 2610
$ws.Range("K116").Value = 21561
$ws.Range("L116").Value = 2610
$ws.Range("M116").Value = -18119
$ws.Range("N116").Value = -9494

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 48488.41
$ws.Range("I132").Value = 61579.117
$ws.Range("J132").Value = 3980
$ws.Range("K132").Value = 184737.351
$ws.Range("L132").Value = 11940
$ws.Range("M132").Value = -182207.351
$ws.Range("N132").Value = -17000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1096.1052
$ws.Range("I137").Value = 894.5
$ws.Range("J137").Value = 1660.6
$ws.Range("K137").Value = 2683.5
$ws.Range("L137").Value = 4981.799999999999
$ws.Range("M137").Value = -133.5
$ws.Range("N137").Value = -10081.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1191.5385
$ws.Range("I61").Value = 1086.3334
$ws.Range("J61").Value = 1428.25
$ws.Range("K61").Value = 1086.3334
$ws.Range("L61").Value = 1428.25
$ws.Range("M61").Value = -874.3334
$ws.Range("N61").Value = -1852.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1105.909
$ws.Range("I74").Value = 1022.8823
$ws.Range("J74").Value = 1194.125
$ws.Range("K74").Value = 1022.8823
$ws.Range("L74").Value = 1194.125
$ws.Range("M74").Value = -148.8823
$ws.Range("N74").Value = -2942.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1105.909
$ws.Range("I77").Value = 1022.8823
$ws.Range("J77").Value = 1194.125
$ws.Range("K77").Value = 5114.4115
$ws.Range("L77").Value = 5970.625
$ws.Range("M77").Value = -746.4115000000002
$ws.Range("N77").Value = -14706.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 587.8387
$ws.Range("I110").Value = 502.46155
$ws.Range("J110").Value = 1031.8
$ws.Range("K110").Value = 502.46155
$ws.Range("L110").Value = 1031.8
$ws.Range("M110").Value = 1542.53845
$ws.Range("N110").Value = -5121.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20343.19
$ws.Range("I132").Value = 1179.439
$ws.Range("J132").Value = 85819.336
$ws.Range("K132").Value = 3538.317
$ws.Range("L132").Value = 257458.008
$ws.Range("M132").Value = -1008.317
$ws.Range("N132").Value = -262518.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1191.5385
$ws.Range("I136").Value = 1086.3334
$ws.Range("J136").Value = 1428.25
$ws.Range("K136").Value = 3259.0002
$ws.Range("L136").Value = 4284.75
$ws.Range("M136").Value = -709.0001999999999
$ws.Range("N136").Value = -9384.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 33020.883
$ws.Range("I86").Value = 2014.45
$ws.Range("J86").Value = 77315.78999999999
$ws.Range("K86").Value = 2014.45
$ws.Range("L86").Value = 77315.78999999999
$ws.Range("M86").Value = -891.45
$ws.Range("N86").Value = -79561.78999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 33020.883
$ws.Range("I89").Value = 2014.45
$ws.Range("J89").Value = 77315.78999999999
$ws.Range("K89").Value = 10072.25
$ws.Range("L89").Value = 386578.95
$ws.Range("M89").Value = -4456.25
$ws.Range("N89").Value = -397810.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4373.5
$ws.Range("I134").Value = 4909.3335
$ws.Range("J134").Value = 2766
$ws.Range("K134").Value = 14728.0005
$ws.Range("L134").Value = 8298
$ws.Range("M134").Value = -12193.0005
$ws.Range("N134").Value = -13368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6479.5454
$ws.Range("I31").Value = 7161.7646
$ws.Range("J31").Value = 4160
$ws.Range("K31").Value = 7161.7646
$ws.Range("L31").Value = 4160
$ws.Range("M31").Value = -6866.7646
$ws.Range("N31").Value = -4750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6479.5454
$ws.Range("I34").Value = 7161.7646
$ws.Range("J34").Value = 4160
$ws.Range("K34").Value = 7161.7646
$ws.Range("L34").Value = 4160
$ws.Range("M34").Value = -6959.7646
$ws.Range("N34").Value = -4564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3672.3333
$ws.Range("I58").Value = 2164
$ws.Range("J58").Value = 3923.7222
$ws.Range("K58").Value = 2164
$ws.Range("L58").Value = 3923.7222
$ws.Range("M58").Value = -1961
$ws.Range("N58").Value = -4329.7222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 114766.555
$ws.Range("I132").Value = 3180
$ws.Range("J132").Value = 254249.75
$ws.Range("K132").Value = 9540
$ws.Range("L132").Value = 762749.25
$ws.Range("M132").Value = -7010
$ws.Range("N132").Value = -767809.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4091.2
$ws.Range("I134").Value = 2073.1428
$ws.Range("J134").Value = 8800
$ws.Range("K134").Value = 6219.428400000001
$ws.Range("L134").Value = 26400
$ws.Range("M134").Value = -3684.428400000001
$ws.Range("N134").Value = -31470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3672.3333
$ws.Range("I136").Value = 2164
$ws.Range("J136").Value = 3923.7222
$ws.Range("K136").Value = 6492
$ws.Range("L136").Value = 11771.1666
$ws.Range("M136").Value = -3942
$ws.Range("N136").Value = -16871.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1794.4445
$ws.Range("I132").Value = 1072.1428
$ws.Range("J132").Value = 2254.0908
$ws.Range("K132").Value = 9649.2852
$ws.Range("L132").Value = 20286.8172
$ws.Range("M132").Value = -7119.2852
$ws.Range("N132").Value = -25346.8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2913.9285
$ws.Range("I80").Value = 1855.7142
$ws.Range("J80").Value = 3972.1428
$ws.Range("K80").Value = 1855.7142
$ws.Range("L80").Value = 3972.1428
$ws.Range("M80").Value = -857.7141999999999
$ws.Range("N80").Value = -5968.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2913.9285
$ws.Range("I83").Value = 1855.7142
$ws.Range("J83").Value = 3972.1428
$ws.Range("K83").Value = 9278.571
$ws.Range("L83").Value = 19860.714
$ws.Range("M83").Value = -4286.571
$ws.Range("N83").Value = -29844.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 52071
$ws.Range("I132").Value = 1431.3846
$ws.Range("J132").Value = 146116
$ws.Range("K132").Value = 4294.1538
$ws.Range("L132").Value = 438348
$ws.Range("M132").Value = -1764.1538
$ws.Range("N132").Value = -443408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1850.7059
$ws.Range("I82").Value = 1905.8
$ws.Range("J82").Value = 1772
$ws.Range("K82").Value = 1905.8
$ws.Range("L82").Value = 1772
$ws.Range("M82").Value = -1544.8
$ws.Range("N82").Value = -2494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1850.7059
$ws.Range("I85").Value = 1905.8
$ws.Range("J85").Value = 1772
$ws.Range("K85").Value = 1905.8
$ws.Range("L85").Value = 1772
$ws.Range("M85").Value = -657.8
$ws.Range("N85").Value = -4268

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 314443.97
$ws.Range("I132").Value = 346386.53
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 1039159.59
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -1036629.59
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 22841.715
$ws.Range("I136").Value = 35401.332
$ws.Range("J136").Value = 13422
$ws.Range("K136").Value = 106203.996
$ws.Range("L136").Value = 40266
$ws.Range("M136").Value = -103653.996
$ws.Range("N136").Value = -45366

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 98637860
$ws.Range("I132").Value = 178572450
$ws.Range("J132").Value = 5380844.5
$ws.Range("K132").Value = 535717350
$ws.Range("L132").Value = 16142533.5
$ws.Range("M132").Value = -535714820
$ws.Range("N132").Value = -16147593.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 23111.09
$ws.Range("I136").Value = 39191.883
$ws.Range("J136").Value = 1105.7894
$ws.Range("K136").Value = 117575.649
$ws.Range("L136").Value = 3317.3682
$ws.Range("M136").Value = -115025.649
$ws.Range("N136").Value = -8417.368200000001
